$wb = $excel.ActiveWorkbook

# --- Rename headers on existing sheets ---
$wsWeekly = $wb.Worksheets.Item("Weekly Quantity")
$wsWeekly.Range("B1").Value = "Weekly_PO_Qty"

$wsMonthly = $wb.Worksheets.Item("Monthly Trend")
$wsMonthly.Range("B1").Value = "Monthly_PO_Qty"

# --- Add new "PO Forecast" sheet after "Monthly Trend" ---
$wsForecast = $wb.Worksheets.Add($null, $wsMonthly)
$wsForecast.Name = "PO Forecast"

# Give the new header row (A1:D1) the same bold/centered/boxed style used
# by the existing sheets' header rows, then fill in the header text
$wsWeekly.Range("A1:B1").Copy() | Out-Null
$wsForecast.Range("A1:D1").PasteSpecial(-4122) | Out-Null  # xlPasteFormats

$wsForecast.Range("A1").Value = "ds"
$wsForecast.Range("B1").Value = "PO_Forecast"
$wsForecast.Range("C1").Value = "yhat_lower"
$wsForecast.Range("D1").Value = "yhat_upper"

# Data rows
$data = @(
    @(45109.99999999999, 29, 27.66894393932442, 31.20935210422522),
    @(45158.99999999999, 22, 20.24219040814769, 23.73864463842471),
    @(45179.99999999999, 19, 16.93067080452278, 20.48014910364848),
    @(45186.99999999999, 18, 15.87693521628075, 19.3249143913227),
    @(45193.99999999999, 17, 14.81638392231285, 18.3295284417951),
    @(45200.99999999999, 15, 13.77199611806492, 17.20011848383992),
    @(45207.99999999999, 14, 12.56385349535434, 16.13965853336696),
    @(45214.99999999999, 13, 11.58034805390781, 15.13160884604513),
    @(45221.99999999999, 12, 10.5884088824303, 13.95568359425618),
    @(45228.99999999999, 11, 9.387978950022616, 12.96969114421188),
    @(45235.99999999999, 10, 8.400824868644111, 11.92239646591486)
)

$row = 2
foreach ($r in $data) {
    $wsForecast.Cells.Item($row, 1).Value = $r[0]
    $wsForecast.Cells.Item($row, 2).Value = $r[1]
    $wsForecast.Cells.Item($row, 3).Value = $r[2]
    $wsForecast.Cells.Item($row, 4).Value = $r[3]
    $row++
}

# Column A ("ds") carries the same date/time number format used by the
# "Order Week" / "Order Month" columns on the other two sheets
$wsForecast.Range("A2:A12").NumberFormat = "YYYY-MM-DD HH:MM:SS"

# Keep the originally-active sheet selected, matching the source workbook
$wsWeekly.Activate() | Out-Null
$wsWeekly.Range("A1").Select() | Out-Null
